# 8.7.1 workbook update
# Mirrors the upstream commit: a handful of label cells on the (only) sheet
# were re-worded (sex/urbanisation/education labels normalised and several
# Kyrgyz/English translations corrected), and a previously-empty header
# cell (A35) was filled in with its missing Kyrgyz label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: "by sex" header row -------------------------------------------------
$ws.Range("A6").Value = "Жынысы боюнча"
$ws.Range("B6").Value = "По полу"
$ws.Range("C6").Value = "By sex"

# --- Row 7: Male ------------------------------------------------------------
$ws.Range("A7").Value = "Эркектер"
$ws.Range("B7").Value = "Мужчины"
$ws.Range("C7").Value = "Men"

# --- Row 8: Female ----------------------------------------------------------
$ws.Range("A8").Value = "Аялдар"
$ws.Range("B8").Value = "Женщины"
$ws.Range("C8").Value = "Woman"

# --- Row 10: urban ------------------------------------------------------------
$ws.Range("A10").Value = "Шаар"
$ws.Range("C10").Value = "Urban"

# --- Row 11: rural ------------------------------------------------------------
$ws.Range("A11").Value = "Айыл"
$ws.Range("C11").Value = "Rural"

# --- School attendance / mother's education (English column) ----------------
$ws.Range("C28").Value = "Does not attend"
$ws.Range("C29").Value = "Educationof mother"
$ws.Range("C30").Value = "Preschool or not /primary"
$ws.Range("C31").Value = "Basic general"
$ws.Range("C32").Value = "Average total"
$ws.Range("C33").Value = "Vocational primary /secondary"
$ws.Range("C34").Value = "Higher"

# --- Row 35: previously-blank "functional difficulties in a child" header ---
# Match the bold-italic styling already used by its row siblings (B35/C35) by
# copying their format over before setting the text.
$ws.Range("B35").Copy() | Out-Null
$ws.Range("A35").PasteSpecial(-4122) | Out-Null
$ws.Range("A35").Value = "Баланын функционалдык кыйнчылыктары"

# --- Wealth quintile (English column) ----------------------------------------
$ws.Range("C38").Value = "Wealth quintile"
